# Generate Report for Handoff
#
# The source markdown file was renamed/regenerated
# (a0584535-...-ca56.md -> 9cdbd6f6-...-d221.md) and two screenshot (.png)
# files that it depends on are now also tracked for localization handoff.
# This updates the Overview sheet and the per-locale (zh-cn / de-de)
# handoff-status sheets to add rows for the two new dependency images and
# to refresh the file names / timestamps for the re-generated handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Shared constants
# ---------------------------------------------------------------------
$newMdName   = "9cdbd6f6-f6e0-47f0-93a4-35aac997d221.md"
$newXlfHash  = "fe7ef8ef85b96ebae73de6f61d0b39d34828b5c0"

$png1 = "c6339246-4958-4aec-ae68-f0672f6c56aa.png"
$png2 = "f9f1385c-ef28-427c-8cb4-f786ddf59ebb.png"
$png1Target = "f0cae7b5404b315d6b3f7b4c65852e5255d74d3c.png"
$png2Target = "166203a53b6772eaa901c2dc910a3b999fdec43e.png"

$configName = ".localization-config"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$includeText     = "Include"
$isDependencyText= "IsDependency"
$ignoredText     = "Ignored"
$zeroDatetime    = "0001-01-01 00:00:00"
$dependencyFromMd = "e2e\9cdbd6f6-f6e0-47f0-93a4-35aac997d221.md"

$zhXlfName = "9cdbd6f6-f6e0-47f0-93a4-35aac997d221.$newXlfHash.zh-cn.xlf"
$deXlfName = "9cdbd6f6-f6e0-47f0-93a4-35aac997d221.$newXlfHash.de-de.xlf"
$zhHandoffDatetime = "2016-03-09 10:04:10"
$deHandoffDatetime = "2016-03-09 10:04:14"

$repoSha       = "a30bf9fceb744902272e9d2d45d7578b559f66e2"
$zhHandoffSha  = "31aea83612b700e827340e7cba4343a97e4152b8"
$deHandoffSha  = "58ed5dcc77d96f367f29c37d74e31306ce9ee196"

function E2EUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$repoSha/e2e/$name"
}
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$repoSha/.localization-config"

function HtZhUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$name"
}
function HtDeUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$name"
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value2 = $newMdName
$ws1.Range("B2").Value2 = $readyForHandoff
$ws1.Range("C2").Value2 = $readyForHandoff

$ws1.Range("A3").Value2 = $png1
$ws1.Range("B3").Value2 = $readyForHandoff
$ws1.Range("C3").Value2 = $readyForHandoff

$ws1.Range("A4").Value2 = $png2
$ws1.Range("B4").Value2 = $readyForHandoff
$ws1.Range("C4").Value2 = $readyForHandoff

$ws1.Range("A5").Value2 = $configName
$ws1.Range("B5").Value2 = $notLocalized
$ws1.Range("C5").Value2 = $notLocalized

$ws1.Hyperlinks.Add($ws1.Range("A2"), (E2EUrl $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), (E2EUrl $png1), [Type]::Missing, [Type]::Missing, $png1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), (E2EUrl $png2), [Type]::Missing, [Type]::Missing, $png2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), $configUrl, [Type]::Missing, [Type]::Missing, $configName) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value2 = $newMdName
$ws2.Range("B2").Value2 = $readyForHandoff
$ws2.Range("C2").Value2 = $zhXlfName
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D2").Value2 = $zhHandoffDatetime
$ws2.Range("G2").Value2 = $zeroDatetime
$ws2.Range("H2").Value2 = $includeText

$ws2.Range("A3").Value2 = $png1
$ws2.Range("B3").Value2 = $readyForHandoff
$ws2.Range("C3").Value2 = $png1Target
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D3").Value2 = $zhHandoffDatetime
$ws2.Range("G3").Value2 = $zeroDatetime
$ws2.Range("H3").Value2 = $isDependencyText
$ws2.Range("I3").Value2 = $dependencyFromMd

$ws2.Range("A4").Value2 = $png2
$ws2.Range("B4").Value2 = $readyForHandoff
$ws2.Range("C4").Value2 = $png2Target
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D4").Value2 = $zhHandoffDatetime
$ws2.Range("G4").Value2 = $zeroDatetime
$ws2.Range("H4").Value2 = $isDependencyText
$ws2.Range("I4").Value2 = $dependencyFromMd

$ws2.Range("A5").Value2 = $configName
$ws2.Range("B5").Value2 = $notLocalized
$ws2.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D5").Value2 = $zeroDatetime
$ws2.Range("G5").Value2 = $zeroDatetime
$ws2.Range("H5").Value2 = $ignoredText

$ws2.Hyperlinks.Add($ws2.Range("A2"), (E2EUrl $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), (HtZhUrl $zhXlfName), [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), (E2EUrl $png1), [Type]::Missing, [Type]::Missing, $png1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), (HtZhUrl $png1Target), [Type]::Missing, [Type]::Missing, $png1Target) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), (E2EUrl $png2), [Type]::Missing, [Type]::Missing, $png2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), (HtZhUrl $png2Target), [Type]::Missing, [Type]::Missing, $png2Target) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), $configUrl, [Type]::Missing, [Type]::Missing, $configName) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value2 = $newMdName
$ws3.Range("B2").Value2 = $readyForHandoff
$ws3.Range("C2").Value2 = $deXlfName
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D2").Value2 = $deHandoffDatetime
$ws3.Range("G2").Value2 = $zeroDatetime
$ws3.Range("H2").Value2 = $includeText

$ws3.Range("A3").Value2 = $png1
$ws3.Range("B3").Value2 = $readyForHandoff
$ws3.Range("C3").Value2 = $png1Target
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D3").Value2 = $deHandoffDatetime
$ws3.Range("G3").Value2 = $zeroDatetime
$ws3.Range("H3").Value2 = $isDependencyText
$ws3.Range("I3").Value2 = $dependencyFromMd

$ws3.Range("A4").Value2 = $png2
$ws3.Range("B4").Value2 = $readyForHandoff
$ws3.Range("C4").Value2 = $png2Target
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D4").Value2 = $deHandoffDatetime
$ws3.Range("G4").Value2 = $zeroDatetime
$ws3.Range("H4").Value2 = $isDependencyText
$ws3.Range("I4").Value2 = $dependencyFromMd

$ws3.Range("A5").Value2 = $configName
$ws3.Range("B5").Value2 = $notLocalized
$ws3.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D5").Value2 = $zeroDatetime
$ws3.Range("G5").Value2 = $zeroDatetime
$ws3.Range("H5").Value2 = $ignoredText

$ws3.Hyperlinks.Add($ws3.Range("A2"), (E2EUrl $newMdName), [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), (HtDeUrl $deXlfName), [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), (E2EUrl $png1), [Type]::Missing, [Type]::Missing, $png1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), (HtDeUrl $png1Target), [Type]::Missing, [Type]::Missing, $png1Target) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), (E2EUrl $png2), [Type]::Missing, [Type]::Missing, $png2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), (HtDeUrl $png2Target), [Type]::Missing, [Type]::Missing, $png2Target) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), $configUrl, [Type]::Missing, [Type]::Missing, $configName) | Out-Null

Write-Host "Done updating localization-status report."
